# "started PA2, added buttons and second layout"
# Fill in the whiteboarding signup roster (column B) for the Spring 2018
# CS 480 schedule, and nudge the view to match the saved state
# (zoomed-in second layout, active selection on B3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$signups = @{
    2  = "Jack Kinne"
    3  = "Zach Freeman"
    4  = "Levi Pole"
    5  = "Chase Cullen"
    6  = "Sam Alston"
    7  = "Jakob Konicke"
    8  = "Lam Ngo"
    9  = "Alec Levin"
    10 = "Jeremy Walker"
    11 = "Jack Witherell"
    12 = "Sam Papavasiliou"
    13 = "Angel Ruiz"
    15 = "Casey Tran"
    16 = "Eric Mott"
    17 = "Julian Lucas"
    18 = "Chadwick Davis"
    19 = "Raymond Ivey"
    20 = "Christopher Johnstone"
    21 = "Marcellus Parley"
    22 = "David Lynch"
    23 = "Luis Garcia"
    24 = "William Dounda"
    25 = "Jared Conn"
    26 = "Thomas Murphy"
}

foreach ($row in $signups.Keys) {
    $ws.Cells.Item($row, 2).Value = $signups[$row]
}

# New working selection + zoomed-in layout for the second (zoomed) view.
[void]$ws.Range("B3").Select()
$excel.ActiveWindow.Zoom = 160
